# Update "想去人数" (interest/want-to-go count) figures in column F
# for a handful of events, reflecting a refreshed data pull.
# Affected sheets: 展览 (Exhibitions) and 全部类型 (All types)

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    if ($name -eq "全部类型") {
        # In this combined sheet, the exhibition rows are interleaved
        # with performance rows, so the 万圣漫控 row lands on F9.
        $ws.Range("F3").Value = 119
        $ws.Range("F4").Value = 153
        $ws.Range("F5").Value = 3085
        $ws.Range("F6").Value = 312
        $ws.Range("F9").Value = 411
    }
    else {
        $ws.Range("F3").Value = 119
        $ws.Range("F4").Value = 153
        $ws.Range("F5").Value = 3085
        $ws.Range("F6").Value = 312
        $ws.Range("F7").Value = 411
    }
}
